{"js": "// The document starts with 10 paragraphs (bulleted list items + a couple of\n// plain paragraphs). The edit:\n//   1. Deletes the first three bullet items (\"Convertir pixeles a mm...\",\n//      \"Analiticamente, calcular...\" and \"Cuanta fuerza se necesita...\").\n//   2. The remaining items shift up. The (now first) item - previously\n//      \"Para el estado tensorial a 1cm...\" - gets its text re-typed with\n//      Word's spell-check markers (<w:proofErr/>) wrapping the unrecognised\n//      tokens s_xx, s_yy and t_xy.\n//   3. The (now fifth) item - \"Buscar e implementar tensor de\n//      deformaciones... para hacer los mapas por mi mismo\" - gets the lone\n//      word \"mi\" split into its own spell-checked run.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// 1. Remove the first three list items; the rest shift up to fill the gap.\nparagraphs.items[0].delete();\nparagraphs.items[1].delete();\nparagraphs.items[2].delete();\nawait context.sync();\n\n// 2. Re-insert the \"estado tensorial\" paragraph (now index 0) with proofErr\n// wrapped spans around s_xx, s_yy and t_xy.\nconst afterDelete = body.paragraphs;\nafterDelete.load(\"items\");\nawait context.sync();\n\nconst tensorialOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"Prrafodelista\"/>\n              <w:numPr>\n                <w:ilvl w:val=\"0\"/>\n                <w:numId w:val=\"1\"/>\n              </w:numPr>\n              <w:jc w:val=\"both\"/>\n            </w:pPr>\n            <w:r>\n              <w:t xml:space=\"preserve\">Para el estado tensorial a 1cm obtener las tensiones que se est\u00e1n aplicando a partir de la matriz que forma los mapas, as\u00ed podemos saber </w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r>\n              <w:t>s_xx</w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r>\n              <w:t xml:space=\"preserve\">, </w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r>\n              <w:t>s_yy</w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r>\n              <w:t xml:space=\"preserve\"> y </w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r>\n              <w:t>t_xy</w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r>\n              <w:t>. Y despu\u00e9s podemos hacer el c\u00edrculo de Mohr.</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\nafterDelete.items[0].insertOoxml(tensorialOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// 3. Re-insert the \"Buscar e implementar tensor...\" paragraph (now index 4)\n// splitting out \"mi\" into its own proofErr wrapped run.\nconst afterFirstFix = body.paragraphs;\nafterFirstFix.load(\"items\");\nawait context.sync();\n\nconst tensorDeformacionesOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"Prrafodelista\"/>\n              <w:numPr>\n                <w:ilvl w:val=\"0\"/>\n                <w:numId w:val=\"1\"/>\n              </w:numPr>\n              <w:jc w:val=\"both\"/>\n            </w:pPr>\n            <w:r>\n              <w:t>Buscar e implementar tensor de deformaciones</w:t>\n            </w:r>\n            <w:r>\n              <w:t xml:space=\"preserve\"> para hacer los mapas por </w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r>\n              <w:t>mi</w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r>\n              <w:t xml:space=\"preserve\"> mismo</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\nafterFirstFix.items[4].insertOoxml(tensorDeformacionesOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The document starts with 10 paragraphs (bulleted list items + a couple of\n# plain paragraphs). The edit:\n#   1. Deletes the first three bullet items (\"Convertir pixeles a mm...\",\n#      \"Analiticamente, calcular...\" and \"Cuanta fuerza se necesita...\").\n#   2. The remaining items shift up. The (now first) item - previously\n#      \"Para el estado tensorial a 1cm...\" - gets its text re-typed with\n#      Word's spell-check markers (proofErr) wrapping the unrecognised\n#      tokens s_xx, s_yy and t_xy.\n#   3. The (now fifth) item - \"Buscar e implementar tensor de\n#      deformaciones... para hacer los mapas por mi mismo\" - gets the lone\n#      word \"mi\" split into its own spell-checked run.\n\n$d = $word.ActiveDocument\n\n# 1. Remove the first three list items; later paragraphs shift up to fill\n# the gap left behind.\n$d.Paragraphs(1).Range.Delete()\n$d.Paragraphs(1).Range.Delete()\n$d.Paragraphs(1).Range.Delete()\n\n# 2. Re-insert the \"estado tensorial\" paragraph (now paragraph 1) with\n# proofErr wrapped spans around s_xx, s_yy and t_xy.\n$tensorialXml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:pStyle w:val=\"Prrafodelista\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr><w:jc w:val=\"both\"/></w:pPr><w:r><w:t xml:space=\"preserve\">Para el estado tensorial a 1cm obtener las tensiones que se est\u00e1n aplicando a partir de la matriz que forma los mapas, as\u00ed podemos saber </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>s_xx</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">, </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>s_yy</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> y </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>t_xy</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>. Y despu\u00e9s podemos hacer el c\u00edrculo de Mohr.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$d.Paragraphs(1).Range.InsertXML($tensorialXml)\n\n# 3. Re-insert the \"Buscar e implementar tensor...\" paragraph (now\n# paragraph 5) splitting out \"mi\" into its own proofErr wrapped run.\n$tensorDeformacionesXml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:pStyle w:val=\"Prrafodelista\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr><w:jc w:val=\"both\"/></w:pPr><w:r><w:t>Buscar e implementar tensor de deformaciones</w:t></w:r><w:r><w:t xml:space=\"preserve\"> para hacer los mapas por </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>mi</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> mismo</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$d.Paragraphs(5).Range.InsertXML($tensorDeformacionesXml)\n"}
